$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.10089999999999
$ws.Range("C21").Value = -13.00690000000001
$ws.Range("C23").Value = -11.98620000000001
$ws.Range("C25").Value = -11.1566
